# Create advisors table student-id label update, clear stray Courses-table
# rows (Semester / Prerequisities), add Campus to the Housing table, and
# give the Advisors Table column the same highlight fill as the other
# table bodies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Advisors Table: rename "Student id" row label to "Students"
$ws.Range("B15").Value = "Students"

# Courses table: remove the Semester / Prerequisities rows (clear content
# and drop their old blue highlight so they go back to a plain cell)
$ws.Range("A22:A23").ClearContents()
$ws.Range("A22:A23").ClearFormats()

# Housing table: add a new "Campus" attribute row
$ws.Range("B25").Value = "Campus"

# Advisors Table body (header + 3 attribute rows) now shares the same
# green fill used by the other table bodies (e.g. Student Table column A/B)
$ws.Range("B13:B16").Interior.Color = 5296274

# Restore the view: scroll back to the top and select I12
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("I12").Select()
